$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update point_id values for several landmark rows; others get cleared.
$ws.Range("B2").Value = 4869
$ws.Range("B3").Value = 3410
$ws.Range("B4").Value = 2431
$ws.Range("B5").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("B10").Value = 9395
$ws.Range("B11").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("B13").Value = 3607
$ws.Range("B14").ClearContents()
$ws.Range("B15").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("B17").ClearContents()
$ws.Range("B18").ClearContents()
$ws.Range("B19").Value = 8948
$ws.Range("B20").Value = 678
$ws.Range("B21").Value = 2102
$ws.Range("B22").Value = 4197
$ws.Range("B23").Value = 2063
$ws.Range("B24").ClearContents()
$ws.Range("B25").ClearContents()
$ws.Range("B26").Value = 8967
$ws.Range("B27").ClearContents()
$ws.Range("B28").ClearContents()

# Update the sheet view: remove frozen/scrolled topLeftCell and change selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("M7").Select()
